$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace AVERAGE(...) with MEDIAN(...) in every summary-row formula cell
# (row blocks 36, 73, 110, 147, 184, 221, 258, 295, 332 across the
# C:F / J:M / Q:T / X:AA column groups). Every cell in each group is set
# explicitly so shared-formula groups don't end up with stale AVERAGE
# formulas left behind.
$updates = @(
    @("C36", "=MEDIAN(C6:C35)"),
    @("D36", "=MEDIAN(D6:D35)"),
    @("E36", "=MEDIAN(E6:E35)"),
    @("F36", "=MEDIAN(F6:F35)"),
    @("J36", "=MEDIAN(J6:J35)"),
    @("K36", "=MEDIAN(K6:K35)"),
    @("L36", "=MEDIAN(L6:L35)"),
    @("M36", "=MEDIAN(M6:M35)"),
    @("Q36", "=MEDIAN(Q6:Q35)"),
    @("R36", "=MEDIAN(R6:R35)"),
    @("S36", "=MEDIAN(S6:S35)"),
    @("T36", "=MEDIAN(T6:T35)"),
    @("X36", "=MEDIAN(X6:X35)"),
    @("Y36", "=MEDIAN(Y6:Y35)"),
    @("Z36", "=MEDIAN(Z6:Z35)"),
    @("AA36", "=MEDIAN(AA6:AA35)"),
    @("C73", "=MEDIAN(C43:C72)"),
    @("D73", "=MEDIAN(D43:D72)"),
    @("E73", "=MEDIAN(E43:E72)"),
    @("F73", "=MEDIAN(F43:F72)"),
    @("J73", "=MEDIAN(J43:J72)"),
    @("K73", "=MEDIAN(K43:K72)"),
    @("L73", "=MEDIAN(L43:L72)"),
    @("M73", "=MEDIAN(M43:M72)"),
    @("Q73", "=MEDIAN(Q43:Q72)"),
    @("R73", "=MEDIAN(R43:R72)"),
    @("S73", "=MEDIAN(S43:S72)"),
    @("T73", "=MEDIAN(T43:T72)"),
    @("X73", "=MEDIAN(X43:X72)"),
    @("Y73", "=MEDIAN(Y43:Y72)"),
    @("Z73", "=MEDIAN(Z43:Z72)"),
    @("AA73", "=MEDIAN(AA43:AA72)"),
    @("C110", "=MEDIAN(C80:C109)"),
    @("D110", "=MEDIAN(D80:D109)"),
    @("E110", "=MEDIAN(E80:E109)"),
    @("F110", "=MEDIAN(F80:F109)"),
    @("J110", "=MEDIAN(J80:J109)"),
    @("K110", "=MEDIAN(K80:K109)"),
    @("L110", "=MEDIAN(L80:L109)"),
    @("M110", "=MEDIAN(M80:M109)"),
    @("Q110", "=MEDIAN(Q80:Q109)"),
    @("R110", "=MEDIAN(R80:R109)"),
    @("S110", "=MEDIAN(S80:S109)"),
    @("T110", "=MEDIAN(T80:T109)"),
    @("X110", "=MEDIAN(X80:X109)"),
    @("Y110", "=MEDIAN(Y80:Y109)"),
    @("Z110", "=MEDIAN(Z80:Z109)"),
    @("AA110", "=MEDIAN(AA80:AA109)"),
    @("C147", "=MEDIAN(C117:C146)"),
    @("D147", "=MEDIAN(D117:D146)"),
    @("E147", "=MEDIAN(E117:E146)"),
    @("F147", "=MEDIAN(F117:F146)"),
    @("J147", "=MEDIAN(J117:J146)"),
    @("K147", "=MEDIAN(K117:K146)"),
    @("L147", "=MEDIAN(L117:L146)"),
    @("M147", "=MEDIAN(M117:M146)"),
    @("Q147", "=MEDIAN(Q117:Q146)"),
    @("R147", "=MEDIAN(R117:R146)"),
    @("S147", "=MEDIAN(S117:S146)"),
    @("T147", "=MEDIAN(T117:T146)"),
    @("X147", "=MEDIAN(X117:X146)"),
    @("Y147", "=MEDIAN(Y117:Y146)"),
    @("Z147", "=MEDIAN(Z117:Z146)"),
    @("AA147", "=MEDIAN(AA117:AA146)"),
    @("C184", "=MEDIAN(C154:C183)"),
    @("D184", "=MEDIAN(D154:D183)"),
    @("E184", "=MEDIAN(E154:E183)"),
    @("F184", "=MEDIAN(F154:F183)"),
    @("J184", "=MEDIAN(J154:J183)"),
    @("K184", "=MEDIAN(K154:K183)"),
    @("L184", "=MEDIAN(L154:L183)"),
    @("M184", "=MEDIAN(M154:M183)"),
    @("Q184", "=MEDIAN(Q154:Q183)"),
    @("R184", "=MEDIAN(R154:R183)"),
    @("S184", "=MEDIAN(S154:S183)"),
    @("T184", "=MEDIAN(T154:T183)"),
    @("X184", "=MEDIAN(X154:X183)"),
    @("Y184", "=MEDIAN(Y154:Y183)"),
    @("Z184", "=MEDIAN(Z154:Z183)"),
    @("AA184", "=MEDIAN(AA154:AA183)"),
    @("C221", "=MEDIAN(C191:C220)"),
    @("D221", "=MEDIAN(D191:D220)"),
    @("E221", "=MEDIAN(E191:E220)"),
    @("F221", "=MEDIAN(F191:F220)"),
    @("J221", "=MEDIAN(J191:J220)"),
    @("K221", "=MEDIAN(K191:K220)"),
    @("L221", "=MEDIAN(L191:L220)"),
    @("M221", "=MEDIAN(M191:M220)"),
    @("Q221", "=MEDIAN(Q191:Q220)"),
    @("R221", "=MEDIAN(R191:R220)"),
    @("S221", "=MEDIAN(S191:S220)"),
    @("T221", "=MEDIAN(T191:T220)"),
    @("X221", "=MEDIAN(X191:X220)"),
    @("Y221", "=MEDIAN(Y191:Y220)"),
    @("Z221", "=MEDIAN(Z191:Z220)"),
    @("AA221", "=MEDIAN(AA191:AA220)"),
    @("C258", "=MEDIAN(C228:C257)"),
    @("D258", "=MEDIAN(D228:D257)"),
    @("E258", "=MEDIAN(E228:E257)"),
    @("F258", "=MEDIAN(F228:F257)"),
    @("J258", "=MEDIAN(J228:J257)"),
    @("K258", "=MEDIAN(K228:K257)"),
    @("L258", "=MEDIAN(L228:L257)"),
    @("M258", "=MEDIAN(M228:M257)"),
    @("Q258", "=MEDIAN(Q228:Q257)"),
    @("R258", "=MEDIAN(R228:R257)"),
    @("S258", "=MEDIAN(S228:S257)"),
    @("T258", "=MEDIAN(T228:T257)"),
    @("X258", "=MEDIAN(X228:X257)"),
    @("Y258", "=MEDIAN(Y228:Y257)"),
    @("Z258", "=MEDIAN(Z228:Z257)"),
    @("AA258", "=MEDIAN(AA228:AA257)"),
    @("C295", "=MEDIAN(C265:C294)"),
    @("D295", "=MEDIAN(D265:D294)"),
    @("E295", "=MEDIAN(E265:E294)"),
    @("F295", "=MEDIAN(F265:F294)"),
    @("J295", "=MEDIAN(J265:J294)"),
    @("K295", "=MEDIAN(K265:K294)"),
    @("L295", "=MEDIAN(L265:L294)"),
    @("M295", "=MEDIAN(M265:M294)"),
    @("Q295", "=MEDIAN(Q265:Q294)"),
    @("R295", "=MEDIAN(R265:R294)"),
    @("S295", "=MEDIAN(S265:S294)"),
    @("T295", "=MEDIAN(T265:T294)"),
    @("X295", "=MEDIAN(X265:X294)"),
    @("Y295", "=MEDIAN(Y265:Y294)"),
    @("Z295", "=MEDIAN(Z265:Z294)"),
    @("AA295", "=MEDIAN(AA265:AA294)"),
    @("C332", "=MEDIAN(C302:C331)"),
    @("D332", "=MEDIAN(D302:D331)"),
    @("E332", "=MEDIAN(E302:E331)"),
    @("F332", "=MEDIAN(F302:F331)"),
    @("J332", "=MEDIAN(J302:J331)"),
    @("K332", "=MEDIAN(K302:K331)"),
    @("L332", "=MEDIAN(L302:L331)"),
    @("M332", "=MEDIAN(M302:M331)"),
    @("Q332", "=MEDIAN(Q302:Q331)"),
    @("R332", "=MEDIAN(R302:R331)"),
    @("S332", "=MEDIAN(S302:S331)"),
    @("T332", "=MEDIAN(T302:T331)"),
    @("X332", "=MEDIAN(X302:X331)"),
    @("Y332", "=MEDIAN(Y302:Y331)"),
    @("Z332", "=MEDIAN(Z302:Z331)"),
    @("AA332", "=MEDIAN(AA302:AA331)"),
)

foreach ($pair in $updates) {
    $cellRef = $pair[0]
    $formula = $pair[1]
    $ws.Range($cellRef).Formula = $formula
}

# Update the saved view state (zoom + selection) to match the edit session.
$win = $excel.ActiveWindow
$win.Zoom = 70
$ws.Range("Y341").Select() | Out-Null
